$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting existing rows 30-39 down to 31-40.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44876
$ws.Range("D30").NumberFormat = $ws.Range("D31").NumberFormat
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100101
$ws.Range("H30").Value = "Berries"
$ws.Range("I30").Value = 100101001
$ws.Range("J30").Value = "Arándano (blue)"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 8000
$ws.Range("O30").Value = 9000
$ws.Range("P30").Value = 8500
$ws.Range("Q30").Value = "$/bandeja 2 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 4250
$ws.Range("T30").Value = 2
